$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.281.55'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.089.90'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.09%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '548.21'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.74'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -11.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.083.23'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.11%  '
$ws.Range('E9').Value = '  -3.79%  '
$ws.Range('E10').Value = '  -5.84%  '
$ws.Range('E11').Value = '  -11.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.470'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.92%  '
$ws.Range('E13').Value = '  -6.28%  '
$ws.Range('E14').Value = '  -8.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.588.66'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.295.84'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.111'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.096.53'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.74'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '488.76'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -12.88%  '
$ws.Range('E21').Value = '  -5.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.720'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.25'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -8.17%  '
$ws.Range('E24').Value = '  -4.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.38'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -9.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.47'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -10.01%  '
$ws.Range('E28').Value = '  -7.64%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.97'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -12.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.56'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.14'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.64%  '
$ws.Range('E33').Value = '  -9.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '57.90'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '512.47'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -10.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.02'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.45%  '
$ws.Range('E37').Value = '  -11.67%  '
$ws.Range('E38').Value = '  -12.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.153.58'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.37%  '
$ws.Range('E40').Value = '  -7.60%  '
$ws.Range('E41').Value = '  -7.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.17'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.76%  '
$ws.Range('E43').Value = '  -12.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.259'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.89%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.24'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.59%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.06'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -10.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '120.92'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.99%  '
$ws.Range('E49').Value = '  -4.57%  '
$ws.Range('E50').Value = '  -9.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.04'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -10.13%  '
